$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (identifier / levelOfDescription / extentAndMedium / notes)
$ws.Range("A2:H2").Font.Name = "Calibri"
$ws.Range("B2").Clear() | Out-Null

$ws.Range("A2").Value = "MCH136-1"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 21M | GRAP COUNT NUMER: NONE"

# Restore frozen header pane + select the newly entered row
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A2:H2").Select() | Out-Null
